# Fruta / hortaliza, semanal
# Insert a new weekly price row at row 4 (shifting the existing rows 4-12
# down to 5-13) and populate it with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 4..12 down to 5..13, opening up a blank row 4.
$ws.Rows("4:4").Insert()

# Fill in the new row 4 with this week's record.
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C4").Value = "Los Lagos"
$ws.Range("D4").Value = 44537
$ws.Range("D4").NumberFormat = $ws.Range("D5").NumberFormat
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101001
$ws.Range("J4").Value = "Arándano (blue)"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 400
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5500
$ws.Range("P4").Value = 5250
$ws.Range("Q4").Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Range("R4").Value = "Región del Maule"
$ws.Range("S4").Value = 3500
$ws.Range("T4").Value = 1.5
